# Update the "startup" sheet query table:
# - CasesTab / SamplesTab / FilesTab / StudyFilesTab query text (col B) and
#   the shared StatQuery text (col C) refreshed to match the new Cypher queries
#   (Neutered Status WHERE clause relocated + weight added for CasesTab,
#   File Type / FileName query rewritten, ORDER BY / LIMIT 100 added everywhere).
# - Row heights adjusted to the new wrapped-text extents.
# - Selection moved from B4 to B3 (and implicit scroll back to top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Row 2: CasesTab
$ws.Range("A2").Value = 'CasesTab'
$ws.Range("B2").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
WHERE demo.neutered_indicator IN ["No"]  
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '''') AS Age,
       coalesce(demo.sex, '''') AS Sex,
       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,
       coalesce(diag.best_response, '''') AS `Response to Treatment`,
       coalesce(co.cohort_description, '''') AS `Cohort`
order by c.case_id asc
limit 100'
$ws.Range("C2").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.neutered_indicator IN ["No"]
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'
$ws.Rows(2).RowHeight = 285

# Row 3: SamplesTab
$ws.Range("A3").Value = 'SamplesTab'
$ws.Range("B3").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.neutered_indicator IN ["No"]  
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed,
        coalesce(diag.disease_term,'''') AS Diagnosis, 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`
order by samp.sample_id asc
limit 100'
$ws.Range("C3").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.neutered_indicator IN ["No"]
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'
$ws.Rows(3).RowHeight = 238.5

# Row 4: FilesTab
$ws.Range("A4").Value = 'FilesTab'
$ws.Range("B4").Value = 'MATCH (f:file)-[*]->(c:case)
MATCH (f)-->(parent)
OPTIONAL MATCH (f)-->(samp:sample)
OPTIONAL MATCH (p:program)<--(s:study)<--(c)
OPTIONAL MATCH (c)<--(demo:demographic)
OPTIONAL MATCH (c)<--(diag:diagnosis)
WITH 
    f, c, parent, samp, p, s, demo, diag
WHERE demo.neutered_indicator IN ["No"] 
WITH 
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    2 as precision
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH 
    DISTINCT f, c, parent, samp, p, s, demo, diag, unit,
    round(factor * value)/factor AS size
RETURN
    coalesce(f.file_name, '''') AS `File Name`,
    coalesce(f.file_format, '''') AS `Format`,
    coalesce(f.file_type, '''') AS `File Type`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
    coalesce(labels(parent)[0], '''') AS `Association`,
    coalesce(f.file_description, '''') AS `Description`,
    coalesce(samp.sample_id, '''') AS `Sample ID`,
    coalesce(c.case_id, '''') AS `Case ID`,
    coalesce(demo.breed,'''') AS Breed ,
    coalesce(diag.disease_term,'''') AS Diagnosis
ORDER BY f.file_name asc
   limit 100'
$ws.Range("C4").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.neutered_indicator IN ["No"]
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'
$ws.Rows(4).RowHeight = 235.5

# Row 5: StudyFilesTab
$ws.Range("A5").Value = 'StudyFilesTab'
$ws.Range("B5").Value = 'MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
WHERE demo.neutered_indicator IN ["No"] 
WITH
        DISTINCT f, c, demo, diag, s,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '''') AS `File Name`,
  coalesce(f.file_type, '''') AS `File Type`,
  coalesce("study", '''') AS `Association`,
  coalesce(f.file_description, '''') AS `Description`,
  coalesce(f.file_format, '''') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
  coalesce(s.clinical_study_designation,'''') AS `Study Code`
  order by ''File Name'' asc
  limit 100'
$ws.Range("C5").Value = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.neutered_indicator IN ["No"]
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`'
$ws.Rows(5).RowHeight = 234

# Move the active selection to B3 and scroll the sheet back to the top.
$ws.Range("B3").Select()
